$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values added column-by-column (E, F, G, H, I, J), each column's
# row 3 value followed by its row 4 value, matching shared-string order.
$ws.Range("E3").Value = "SDS"
$ws.Range("E4").Value = "FF"

$ws.Range("F3").Value = "0154545"
$ws.Range("F4").Value = "2445"

$ws.Range("G3").Value = "das"
$ws.Range("G4").Value = "afaf"

$ws.Range("H3").Value = "014545"
$ws.Range("H4").Value = "04554"

$ws.Range("I3").Value = "ad"
$ws.Range("I4").Value = "ddaff"

$ws.Range("J3").Value = "01542"
$ws.Range("J4").Value = "05442"

# Update the selection to match the new active cell
$ws.Range("I9").Select()
